$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": a new date column ("08-dec") is inserted right before
# the old "01-oct." column (EE), shifting EE:FI -> EF:FJ. The freshly
# inserted column gets the new header in row 1 and "-" placeholders for
# every data row (2-25), matching the pattern used by the other untraded
# date columns on this sheet.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Columns("EE:EE").Insert()

$wsSpot.Range("EE1").Value = "08-dec"

for ($row = 2; $row -le 25; $row++) {
    $wsSpot.Cells.Item($row, 135).Value = "-"
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": two new trailing rows for 2025-12-06 / 2025-12-07.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A163").NumberFormat = "@"
$wsGaz.Range("A163").Value = "2025-12-06"
$wsGaz.Range("A163").ClearFormats()
$wsGaz.Range("B163").Value = 25.905

$wsGaz.Range("A164").NumberFormat = "@"
$wsGaz.Range("A164").Value = "2025-12-07"
$wsGaz.Range("A164").ClearFormats()
$wsGaz.Range("B164").Value = 25.905

# ---------------------------------------------------------------------------
# Sheet "CO2": two new trailing rows for 2025-12-06 / 2025-12-07.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A164").NumberFormat = "@"
$wsCo2.Range("A164").Value = "2025-12-06"
$wsCo2.Range("A164").ClearFormats()
$wsCo2.Range("B164").Value = 81.78

$wsCo2.Range("A165").NumberFormat = "@"
$wsCo2.Range("A165").Value = "2025-12-07"
$wsCo2.Range("A165").ClearFormats()
$wsCo2.Range("B165").Value = 81.78
